# Updated cryptos list on Sat Dec 23 19:36:15 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text so that
# numeric-looking strings (e.g. '103.00') are not coerced into numbers
# and keep their original formatting (trailing zeros, padding spaces, etc).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '43.763.36'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '2.295.01'
$ws.Range("E3").Value = '  -1.07%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '103.00'
$ws.Range("E5").Value = '  +5.73%  '
$ws.Range("D6").Value = '270.34'
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  -2.43%  '
$ws.Range("D10").Value = '45.43'
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("D11").Value = '0.0936'
$ws.Range("E11").Value = '  -1.45%  '
$ws.Range("D12").Value = '7.95'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '15.79'
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").Value = '0.857'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '2.302.73'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '43.743.88'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("E19").Value = '  -2.37%  '
$ws.Range("D20").Value = '72.33'
$ws.Range("E20").Value = '  -0.55%  '
$ws.Range("E21").Value = '  +9.33%  '
$ws.Range("D22").Value = '233.41'
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = '2.87'
$ws.Range("E23").Value = '  +14.02%  '
$ws.Range("D24").Value = '9.16'
$ws.Range("E24").Value = '  -2.34%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '11.21'
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("E27").Value = '  -0.64%  '
$ws.Range("D28").Value = '39.92'
$ws.Range("E28").Value = '  +4.85%  '
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("D30").Value = '177.42'
$ws.Range("E30").Value = '  +1.20%  '
$ws.Range("D31").Value = '21.83'
$ws.Range("E31").Value = '  -2.70%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = '5.45'
$ws.Range("E33").Value = '  -0.61%  '
$ws.Range("D34").Value = '4.85'
$ws.Range("E34").Value = '  +10.51%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("D37").Value = '0.0354'
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("E38").Value = '  +6.25%  '
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("D40").Value = '0.235'
$ws.Range("E40").Value = '  -4.10%  '
$ws.Range("D41").Value = '1.40'
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("D42").Value = '12.30'
$ws.Range("E42").Value = '  +1.28%  '
$ws.Range("D43").Value = '65.40'
$ws.Range("E43").Value = '  +5.52%  '
$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").Value = '5.25'
$ws.Range("E44").Value = '  -2.10%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '8.78'
$ws.Range("E45").Value = '  -4.23%  '
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("D47").Value = '1.22'
$ws.Range("E47").Value = '  +1.47%  '
$ws.Range("D48").Value = '98.51'
$ws.Range("E48").Value = '  -1.76%  '
$ws.Range("D49").Value = '0.449'
$ws.Range("E49").Value = '  +8.01%  '
$ws.Range("E50").Value = '  +11.44%  '
$ws.Range("D51").Value = '2.522.74'
$ws.Range("E51").Value = '  -0.80%  '

# Remove the temporary text-number-format marker again so the cells'
# styling matches the original (unstyled) cells.
$ws.Range("D2:E51").ClearFormats()

